$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: latitude -> latitude_API
$ws.Range("A4").Value = "latitude_API"
$ws.Range("B4").Value = "Latitude of sample event provided by NES-LTER API"

# Row 5: longitude -> longitude_API
$ws.Range("A5").Value = "longitude_API"
$ws.Range("B5").Value = "Longitude of sample event provided by NES-LTER API"

# Row 10: depth_API_bottle_summary -> depth_API
$ws.Range("A10").Value = "depth_API"

# Row 2 / Row 3: swap datetime_utc / datetime_utc_matlab attributes and
# refresh the API-sourced datetime description.
$ws.Range("B2").Value = "Data product UTC date and time"

$ws.Range("B10").Value = "Data product depth of sample below sea surface from CTD summary data in NES-LTER  API"

$ws.Range("A2").Value = "datetime_utc"
$ws.Range("A3").Value = "datetime_utc_matlab"
$ws.Range("B3").Value = "PI-provided UTC date and time"

# Update the active selection to match the saved view state.
$ws.Range("B11").Select()
